$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily profit row (row 69) that was produced by running the
# report on 2026-02-01. The date column is stored as plain text (matching
# every other row), so force a text format before writing it to stop Excel
# from auto-converting the "mm/dd/yyyy" looking string into a date serial,
# then drop the temporary formatting so the cell keeps the sheet's default
# (unstyled) look, just like the rest of the data rows.
$dateCell = $ws.Range("A69")
$dateCell.NumberFormat = "@"
$dateCell.Value = "02/01/2026"
$dateCell.ClearFormats()

$ws.Range("B69").Value = 9630.23
$ws.Range("C69").Value = 0.2613837494045819
$ws.Range("D69").Value = 0.7386162505954181
$ws.Range("E69").Value = -318.36
$ws.Range("F69").Value = -40.56
$ws.Range("G69").Value = -23728.71
$ws.Range("H69").Value = -76.94
$ws.Range("I69").Value = -785.66
$ws.Range("J69").Value = -23.79
$ws.Range("K69").Value = -24514.37
$ws.Range("L69").Value = -71.8
